# Remove the row for ADL6172022 (Client.ADLID) from Sheet1.
# This row's Client.Full Name cell had erroneously been pointing at the
# shared string "Mrs Toni Pearce"; the fix is to drop the whole row, which
# shifts every following row up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(288).Delete()
